$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vessel density (tumor)")

# Remove the row referencing "Koyama et al., 2017" (row 2), shifting the
# remaining rows up. This also shrinks the backing table (Table4) from
# A1:C9 to A1:C8 and drops the now-unused shared string.
$ws.Rows.Item(2).Delete()

# Make "Vessel density (tumor)" the active sheet/tab and select the
# (now) second data row, matching the reviewer's final cursor position.
$ws.Activate()
$ws.Range("A2:C2").Select()
